# ---------------------------------------------------------------------------
# Test.xlsx edit: "process conversion of Timestamp and Duration, test nested
# message, list, and map"
#
# Item sheet  (sheet1): add per-attribute "strengthen" nested-message columns
#                        (ID/desc/hint-list/type/count) for attribute 1 and
#                        attribute 2, plus a trailing "valid duration" column.
# Activity sheet (sheet2): add a "duration" column, and split activity
#                        100001 into two chapters/two sections (new rows).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Item
$ws2 = $wb.Worksheets.Item(2)   # Activity

# ===========================================================================
# Sheet1 "Item": insert two blocks of 5 new columns (after attribute-1-Value,
# and after attribute-2-Value), then append one more trailing column.
# ===========================================================================

# Old layout: A..J = 道具ID,名字,描述,IconID,属性1ID,属性1Value,属性2ID,属性2Value,效果列表,过期时间
# Insert 5 blank columns before old G (属性2ID) -> old G,H become L,M
$ws1.Range("G1:K1").EntireColumn.Insert()
# Insert 5 more blank columns before new N -> old I,J (now at S po... ) shift to S,T
$ws1.Range("N1:R1").EntireColumn.Insert()
# New layout is now A..T (20 cols); append one more column U below.

# --- Row 1: headers -------------------------------------------------------
$ws1.Range("A1").Value = "道具ID"
$ws1.Range("B1").Value = "名字"
$ws1.Range("C1").Value = "描述"
$ws1.Range("D1").Value = "IconID"
$ws1.Range("E1").Value = "属性1ID"
$ws1.Range("F1").Value = "属性1Value"
$ws1.Range("G1").Value = "属性1强化1ID"
$ws1.Range("H1").Value = "属性1强化1描述"
$ws1.Range("I1").Value = "属性1强化1提示列表"
$ws1.Range("J1").Value = "属性1类型"
$ws1.Range("K1").Value = "属性1数量"
$ws1.Range("L1").Value = "属性2ID"
$ws1.Range("M1").Value = "属性2Value"
$ws1.Range("N1").Value = "属性2强化1ID"
$ws1.Range("O1").Value = "属性2强化1描述"
$ws1.Range("P1").Value = "属性2强化1提示列表"
$ws1.Range("Q1").Value = "属性2类型"
$ws1.Range("R1").Value = "属性2数量"
$ws1.Range("S1").Value = "效果列表"
$ws1.Range("T1").Value = "过期时间"
$ws1.Range("U1").Value = "有效多久"

# --- Row 2 (item 1001, 金币) ------------------------------------------------
$ws1.Range("A2").Value = 1001
$ws1.Range("B2").Value = "金币"
$ws1.Range("C2").Value = "游戏内通用货币"
$ws1.Range("D2").Value = 5001001
$ws1.Range("E2").Value = 2001
$ws1.Range("F2").Value = 1
$ws1.Range("G2").Value = 1
$ws1.Range("H2").Value = "强化1描述1"
$ws1.Range("I2").Value = "新"
$ws1.Range("J2").Value = 1
$ws1.Range("K2").Value = 1
$ws1.Range("L2").Value = 2002
$ws1.Range("M2").Value = 2
$ws1.Range("N2").Value = 11
$ws1.Range("O2").Value = "强化2描述1"
$ws1.Range("P2").Value = "新"
$ws1.Range("Q2").Value = 1
$ws1.Range("R2").Value = 5
$ws1.Range("S2").Value = 1
$ws1.Range("T2").Value = "2020-01-01  05:00:00"
$ws1.Range("U2").Value = "100"

# --- Row 3 (item 1002, 点券) ------------------------------------------------
$ws1.Range("A3").Value = 1002
$ws1.Range("B3").Value = "点券"
$ws1.Range("C3").Value = "游戏内支付货币"
$ws1.Range("D3").Value = 5001002
$ws1.Range("E3").Value = 2001
$ws1.Range("F3").Value = 1
$ws1.Range("G3").Value = 2
$ws1.Range("H3").Value = "强化1描述2"
$ws1.Range("I3").Value = "新,热"
$ws1.Range("J3").Value = 2
$ws1.Range("K3").Value = 1
$ws1.Range("L3").Value = 2002
$ws1.Range("M3").Value = 2
$ws1.Range("N3").Value = 12
$ws1.Range("O3").Value = "强化2描述2"
$ws1.Range("P3").Value = "新,热"
$ws1.Range("Q3").Value = 2
$ws1.Range("R3").Value = 6
$ws1.Range("S3").Value = "1,2,3"
$ws1.Range("T3").Value = "2020-01-01  05:00:00"
$ws1.Range("U3").Value = "100"

# --- Row 4 (item 1003, 宝石) ------------------------------------------------
$ws1.Range("A4").Value = 1003
$ws1.Range("B4").Value = "宝石"
$ws1.Range("C4").Value = "游戏内稀有货币"
$ws1.Range("D4").Value = 5001003
$ws1.Range("E4").Value = 2001
$ws1.Range("F4").Value = 1
$ws1.Range("G4").Value = 3
$ws1.Range("H4").Value = "强化1描述3"
$ws1.Range("I4").Value = "新,热,限"
$ws1.Range("J4").Value = 3
$ws1.Range("K4").Value = 1
$ws1.Range("L4").Value = 2002
$ws1.Range("M4").Value = 2
$ws1.Range("N4").Value = 13
$ws1.Range("O4").Value = "强化2描述3"
$ws1.Range("P4").Value = "新,热,限"
$ws1.Range("Q4").Value = 3
$ws1.Range("R4").Value = 7
$ws1.Range("S4").Value = "1,2,3,4,5"
$ws1.Range("T4").Value = "2020-01-01  05:00:00"
$ws1.Range("U4").Value = "100"

# Item sheet becomes the active / selected sheet+cell.
$ws1.Activate()
$ws1.Range("P11").Select()

# ===========================================================================
# Sheet2 "Activity": append a trailing "持续多久" column, and split activity
# 100001 (chapter 1) into two sections, plus add a second chapter.
# ===========================================================================

# Insert two new blank rows right after row 2 (100001 / chapter1 / section1),
# pushing the old 100002 / 100003 rows down to rows 5 / 6.
$ws2.Range("A3:A4").EntireRow.Insert()

# --- Row 1: headers --------------------------------------------------------
$ws2.Range("A1").Value = "活动ID"
$ws2.Range("B1").Value = "章ID"
$ws2.Range("C1").Value = "章描述"
$ws2.Range("D1").Value = "节ID"
$ws2.Range("E1").Value = "节描述"
$ws2.Range("F1").Value = "奖励1ID"
$ws2.Range("G1").Value = "奖励1NUM"
$ws2.Range("H1").Value = "奖励2ID"
$ws2.Range("I1").Value = "奖励2NUM"
$ws2.Range("J1").Value = "开始时间"
$ws2.Range("K1").Value = "结束时间"
$ws2.Range("L1").Value = "持续多久"

# --- Row 2: activity 100001, chapter 1, section 1 --------------------------
$ws2.Range("A2").Value = 100001
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = "签到活动章1"
$ws2.Range("D2").Value = 1
$ws2.Range("E2").Value = "签到活动节1"
$ws2.Range("F2").Value = 2001
$ws2.Range("G2").Value = 1
$ws2.Range("H2").Value = 2002
$ws2.Range("I2").Value = 2
$ws2.Range("J2").Value = "2020-01-01  05:00:00"
$ws2.Range("K2").Value = "2020-10-01  05:00:00"
$ws2.Range("L2").Value = "100"

# --- Row 3 (new): activity 100001, chapter 1, section 2 --------------------
$ws2.Range("A3").Value = 100001
$ws2.Range("B3").Value = 1
$ws2.Range("C3").Value = "签到活动章1"
$ws2.Range("D3").Value = 2
$ws2.Range("E3").Value = "签到活动节2"
$ws2.Range("F3").Value = 2001
$ws2.Range("G3").Value = 1
$ws2.Range("H3").Value = 2002
$ws2.Range("I3").Value = 2
$ws2.Range("J3").Value = "2020-01-01  05:00:00"
$ws2.Range("K3").Value = "2020-10-01  05:00:00"
$ws2.Range("L3").Value = "100"

# --- Row 4 (new): activity 100001, chapter 2, section 1 --------------------
$ws2.Range("A4").Value = 100001
$ws2.Range("B4").Value = 2
$ws2.Range("C4").Value = "签到活动章2"
$ws2.Range("D4").Value = 1
$ws2.Range("E4").Value = "签到活动节1"
$ws2.Range("F4").Value = 2002
$ws2.Range("G4").Value = 2
$ws2.Range("H4").Value = 2002
$ws2.Range("I4").Value = 3
$ws2.Range("J4").Value = "2020-01-01  05:00:00"
$ws2.Range("K4").Value = "2020-10-01  05:00:00"
$ws2.Range("L4").Value = "100"

# --- Row 5: activity 100002 (unchanged, shifted down from old row 3) ------
$ws2.Range("A5").Value = 100002
$ws2.Range("B5").Value = 1
$ws2.Range("C5").Value = "抽奖活动章1"
$ws2.Range("D5").Value = 1
$ws2.Range("E5").Value = "抽奖活动节1"
$ws2.Range("F5").Value = 2001
$ws2.Range("G5").Value = 1
$ws2.Range("H5").Value = 2002
$ws2.Range("I5").Value = 2
$ws2.Range("J5").Value = "2020-01-01  05:00:00"
$ws2.Range("K5").Value = "2020-10-01  05:00:00"
$ws2.Range("L5").Value = "100"

# --- Row 6: activity 100003 (unchanged, shifted down from old row 4) ------
$ws2.Range("A6").Value = 100003
$ws2.Range("B6").Value = 1
$ws2.Range("C6").Value = "月卡活动章1"
$ws2.Range("D6").Value = 1
$ws2.Range("E6").Value = "月卡活动节1"
$ws2.Range("F6").Value = 2001
$ws2.Range("G6").Value = 1
$ws2.Range("H6").Value = 2002
$ws2.Range("I6").Value = 2
$ws2.Range("J6").Value = "2020-01-01  05:00:00"
$ws2.Range("K6").Value = "2020-10-01  05:00:00"
$ws2.Range("L6").Value = "100"

# Activity sheet keeps the cursor at its new selection, but Item stays the
# active tab (set above via $ws1.Activate() after this sheet's edits run).
$ws2.Range("J10").Select()

# Re-activate Item last so it is the tab that ends up selected in the file.
$ws1.Activate()
$ws1.Range("P11").Select()
